$d = $word.ActiveDocument

function Add-ReportParagraph([string]$text, [int]$align) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    if ($text -ne "") {
        $r.Text = $text
    }
    $p = $d.Paragraphs.Last
    $p.Alignment = $align
}

# wdAlignParagraphCenter = 1, wdAlignParagraphLeft = 0
Add-ReportParagraph "Báo cáo tuần 9" 1
Add-ReportParagraph "Hoàn thành lại chức năng đăng nhập và đăng ký" 0
Add-ReportParagraph "" 0
Add-ReportParagraph "Báo cáo tuần 10" 1
Add-ReportParagraph "Hoàn thành chức năng CRUD category" 0
